$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 24.11008933333333
$ws.Range("H2").Value = 72.33026799999999
$ws.Range("I2").Value = 0.2144059681009565
$ws.Range("J2").Value = 0.2144059681009565
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.656454333333333
$ws.Range("N2").Value = 4.969363
$ws.Range("O2").Value = 0.0006990759915034363
$ws.Range("P2").Value = 0.0006990759915034364
$ws.Range("Q2").Value = 39.93726195325377
$ws.Range("R2").Value = 359.435357579284
$ws.Range("S2").Value = 0.0001498860647344303
$ws.Range("T2").Value = 0.0001498860647344303

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 24.11008933333333
$ws.Range("H3").Value = 72.33026799999999
$ws.Range("I3").Value = 0.2144059681009565
$ws.Range("J3").Value = 0.2144059681009565
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 6.599386
$ws.Range("N3").Value = 19.798158
$ws.Range("O3").Value = 0.002785149109411345
$ws.Range("P3").Value = 0.002785149109411345
$ws.Range("Q3").Value = 159.1117860051493
$ws.Range("R3").Value = 1432.006074046344
$ws.Range("S3").Value = 0.0005971525911088563
$ws.Range("T3").Value = 0.0005971525911088563

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 24.11008933333333
$ws.Range("H4").Value = 72.33026799999999
$ws.Range("I4").Value = 0.2144059681009565
$ws.Range("J4").Value = 0.2144059681009565
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2347.576497333333
$ws.Range("N4").Value = 7042.729492
$ws.Range("O4").Value = 0.9907513503260664
$ws.Range("P4").Value = 0.9907513503260664
$ws.Range("Q4").Value = 56600.27906754042
$ws.Range("R4").Value = 509402.5116078638
$ws.Range("S4").Value = 0.2124230024139902
$ws.Range("T4").Value = 0.2124230024139902

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 24.11008933333333
$ws.Range("H5").Value = 72.33026799999999
$ws.Range("I5").Value = 0.2144059681009565
$ws.Range("J5").Value = 0.2144059681009565
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 13.65875266666666
$ws.Range("N5").Value = 40.97625799999999
$ws.Range("O5").Value = 0.005764424573018837
$ws.Range("P5").Value = 0.005764424573018838
$ws.Range("Q5").Value = 329.3137469752381
$ws.Range("R5").Value = 2963.823722777143
$ws.Range("S5").Value = 0.001235927031123047
$ws.Range("T5").Value = 0.001235927031123047

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 27.64911833333333
$ws.Range("H6").Value = 82.947355
$ws.Range("I6").Value = 0.2458778107968398
$ws.Range("J6").Value = 0.2458778107968398
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.656454333333333
$ws.Range("N6").Value = 4.969363
$ws.Range("O6").Value = 0.0006990759915034363
$ws.Range("P6").Value = 0.0006990759915034364
$ws.Range("Q6").Value = 45.79950187609612
$ws.Range("R6").Value = 412.1955168848651
$ws.Range("S6").Value = 0.0001718872743714951
$ws.Range("T6").Value = 0.0001718872743714951

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 27.64911833333333
$ws.Range("H7").Value = 82.947355
$ws.Range("I7").Value = 0.2458778107968398
$ws.Range("J7").Value = 0.2458778107968398
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 6.599386
$ws.Range("N7").Value = 19.798158
$ws.Range("O7").Value = 0.002785149109411345
$ws.Range("P7").Value = 0.002785149109411345
$ws.Range("Q7").Value = 182.4672044413433
$ws.Range("R7").Value = 1642.20483997209
$ws.Range("S7").Value = 0.0006848063657648296
$ws.Range("T7").Value = 0.0006848063657648297

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 27.64911833333333
$ws.Range("H8").Value = 82.947355
$ws.Range("I8").Value = 0.2458778107968398
$ws.Range("J8").Value = 0.2458778107968398
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2347.576497333333
$ws.Range("N8").Value = 7042.729492
$ws.Range("O8").Value = 0.9907513503260664
$ws.Range("P8").Value = 0.9907513503260664
$ws.Range("Q8").Value = 64908.42037132152
$ws.Range("R8").Value = 584175.7833418937
$ws.Range("S8").Value = 0.2436037730621861
$ws.Range("T8").Value = 0.2436037730621861

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 27.64911833333333
$ws.Range("H9").Value = 82.947355
$ws.Range("I9").Value = 0.2458778107968398
$ws.Range("J9").Value = 0.2458778107968398
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 13.65875266666666
$ws.Range("N9").Value = 40.97625799999999
$ws.Range("O9").Value = 0.005764424573018837
$ws.Range("P9").Value = 0.005764424573018838
$ws.Range("Q9").Value = 377.6524687663988
$ws.Range("R9").Value = 3398.872218897589
$ws.Range("S9").Value = 0.00141734409451738
$ws.Range("T9").Value = 0.00141734409451738

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 24.21819066666667
$ws.Range("H10").Value = 72.654572
$ws.Range("I10").Value = 0.2153672905874018
$ws.Range("J10").Value = 0.2153672905874018
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.656454333333333
$ws.Range("N10").Value = 4.969363
$ws.Range("O10").Value = 0.0006990759915034363
$ws.Range("P10").Value = 0.0006990759915034364
$ws.Range("Q10").Value = 40.1163268752929
$ws.Range("R10").Value = 361.046941877636
$ws.Range("S10").Value = 0.0001505581022047966
$ws.Range("T10").Value = 0.0001505581022047966

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 24.21819066666667
$ws.Range("H11").Value = 72.654572
$ws.Range("I11").Value = 0.2153672905874018
$ws.Range("J11").Value = 0.2153672905874018
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 6.599386
$ws.Range("N11").Value = 19.798158
$ws.Range("O11").Value = 0.002785149109411345
$ws.Range("P11").Value = 0.002785149109411345
$ws.Range("Q11").Value = 159.8251884309307
$ws.Range("R11").Value = 1438.426695878376
$ws.Range("S11").Value = 0.0005998300175758365
$ws.Range("T11").Value = 0.0005998300175758365

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 24.21819066666667
$ws.Range("H12").Value = 72.654572
$ws.Range("I12").Value = 0.2153672905874018
$ws.Range("J12").Value = 0.2153672905874018
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 2347.576497333333
$ws.Range("N12").Value = 7042.729492
$ws.Range("O12").Value = 0.9907513503260664
$ws.Range("P12").Value = 0.9907513503260664
$ws.Range("Q12").Value = 56854.05521700416
$ws.Range("R12").Value = 511686.4969530375
$ws.Range("S12").Value = 0.2133754339655347
$ws.Range("T12").Value = 0.2133754339655347

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 24.21819066666667
$ws.Range("H13").Value = 72.654572
$ws.Range("I13").Value = 0.2153672905874018
$ws.Range("J13").Value = 0.2153672905874018
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 13.65875266666666
$ws.Range("N13").Value = 40.97625799999999
$ws.Range("O13").Value = 0.005764424573018837
$ws.Range("P13").Value = 0.005764424573018838
$ws.Range("Q13").Value = 330.7902763501751
$ws.Range("R13").Value = 2977.112487151576
$ws.Range("S13").Value = 0.001241468502086507
$ws.Range("T13").Value = 0.001241468502086508

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 36.47324633333334
$ws.Range("H14").Value = 109.419739
$ws.Range("I14").Value = 0.3243489305148018
$ws.Range("J14").Value = 0.3243489305148018
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.656454333333333
$ws.Range("N14").Value = 4.969363
$ws.Range("O14").Value = 0.0006990759915034363
$ws.Range("P14").Value = 0.0006990759915034364
$ws.Range("Q14").Value = 60.41626693958412
$ws.Range("R14").Value = 543.7464024562571
$ws.Range("S14").Value = 0.0002267445501927142
$ws.Range("T14").Value = 0.0002267445501927142

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 36.47324633333334
$ws.Range("H15").Value = 109.419739
$ws.Range("I15").Value = 0.3243489305148018
$ws.Range("J15").Value = 0.3243489305148018
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 6.599386
$ws.Range("N15").Value = 19.798158
$ws.Range("O15").Value = 0.002785149109411345
$ws.Range("P15").Value = 0.002785149109411345
$ws.Range("Q15").Value = 240.7010312267514
$ws.Range("R15").Value = 2166.309281040762
$ws.Range("S15").Value = 0.0009033601349618224
$ws.Range("T15").Value = 0.0009033601349618225

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 36.47324633333334
$ws.Range("H16").Value = 109.419739
$ws.Range("I16").Value = 0.3243489305148018
$ws.Range("J16").Value = 0.3243489305148018
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 2347.576497333333
$ws.Range("N16").Value = 7042.729492
$ws.Range("O16").Value = 0.9907513503260664
$ws.Range("P16").Value = 0.9907513503260664
$ws.Range("Q16").Value = 85623.73587358253
$ws.Range("R16").Value = 770613.6228622426
$ws.Range("S16").Value = 0.3213491408843553
$ws.Range("T16").Value = 0.3213491408843553

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 36.47324633333334
$ws.Range("H17").Value = 109.419739
$ws.Range("I17").Value = 0.3243489305148018
$ws.Range("J17").Value = 0.3243489305148018
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 13.65875266666666
$ws.Range("N17").Value = 40.97625799999999
$ws.Range("O17").Value = 0.005764424573018837
$ws.Range("P17").Value = 0.005764424573018838
$ws.Range("Q17").Value = 498.1790506174069
$ws.Range("R17").Value = 4483.611455556662
$ws.Range("S17").Value = 0.001869684945291903
$ws.Range("T17").Value = 0.001869684945291903

Write-Output "applied 224 cell updates"